# Updates the "Estado de Cuenta" workbook:
#  - Refresh the summary figures (Valor Mora total, worker/period counts)
#  - Replace the worker/period data table with the new single entry
#  - Remove the now-obsolete extra data rows
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block -------------------------------------------------------
$ws.Range("E11").Value = 56940
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# --- Single remaining worker/period row (row 16) -------------------------
$ws.Range("C16").Value = "1007978371"
$ws.Range("D16").Value = "VALENTINA MARTINEZ LEAL"
$ws.Range("E16").Value = "2508"
$ws.Range("G16").Value = 1423500

# --- Drop the rows for the workers/periods that no longer apply ----------
# (rows 17-21 held the other two workers' second period + duplicate entries)
$ws.Rows("17:21").Delete()
